$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "2025/12/04 13:00"
$ws.Range("B44").Value = "-"
$ws.Range("C44").Value = "-"
$ws.Range("D44").Value = "-"
$ws.Range("E44").Value = "-"
$ws.Range("F44").Value = "-"
$ws.Range("G44").Value = "-"
